# Update cryptocurrency price (D) and 1h volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.989.62"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.318.98"
$ws.Range("E3").Value = "  +2.60%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.92"
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.641"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.29"
$ws.Range("E7").Value = "  +7.34%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.652"
$ws.Range("E9").Value = "  -3.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.76"
$ws.Range("E10").Value = "  +1.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0982"
$ws.Range("E11").Value = "  +0.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.59"
$ws.Range("E12").Value = "  -0.30%  "

$ws.Range("E13").Value = "  +2.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.667.70"
$ws.Range("E14").Value = "  +2.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.42"
$ws.Range("E15").Value = "  +3.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.884"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.329.52"
$ws.Range("E17").Value = "  +2.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.007.52"
$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("E19").Value = "  +2.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.35"
$ws.Range("E20").Value = "  +1.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.05"
$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.78"
$ws.Range("E22").Value = "  +1.03%  "

$ws.Range("E23").Value = "  +5.76%  "

$ws.Range("E24").Value = "  -0.78%  "

$ws.Range("E25").Value = "  -0.90%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.43"
$ws.Range("E27").Value = "  -1.21%  "

$ws.Range("E28").Value = "  -1.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "21.31"
$ws.Range("E29").Value = "  +0.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.47"
$ws.Range("E30").Value = "  -0.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0845"
$ws.Range("E31").Value = "  +10.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.28"
$ws.Range("E32").Value = "  +0.43%  "

$ws.Range("E33").Value = "  +1.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.59"
$ws.Range("E34").Value = "  +5.81%  "

$ws.Range("E35").Value = "  +1.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.59"
$ws.Range("E36").Value = "  +11.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.88"
$ws.Range("E37").Value = "  +3.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0315"
$ws.Range("E38").Value = "  -2.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.28"
$ws.Range("E39").Value = "  +17.45%  "

$ws.Range("E40").Value = "  +2.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.92"
$ws.Range("E41").Value = "  +1.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.219"
$ws.Range("E42").Value = "  +8.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.79"
$ws.Range("E43").Value = "  -2.33%  "

$ws.Range("E44").Value = "  +3.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.90"
$ws.Range("E45").Value = "  -2.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.64"
$ws.Range("E46").Value = "  +12.29%  "

$ws.Range("E47").Value = "  -0.58%  "

$ws.Range("E48").Value = "  -0.30%  "

$ws.Range("E49").Value = "  -0.06%  "

$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.35"
$ws.Range("E51").Value = "  -0.64%  "
